# New response model class added — workbook now stores string-typed
# response values ("ghdbjkm" / "cgfvhbjnk") in testData!D4:D5 instead of
# the old numeric placeholders, and the active sheet/selection moves
# from "RestAssured" back to "testData".

$wb = $excel.ActiveWorkbook

$testData = $wb.Worksheets.Item("testData")

# Replace the numeric placeholder values with the new string values.
$testData.Range("D4").Value = "ghdbjkm"
$testData.Range("D5").Value = "cgfvhbjnk"

# Move the active sheet/tab back to "testData" and set its selection.
$testData.Activate()
[void]$testData.Range("J8").Select()
